$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Sheet1"

# Insert 4 new rows at the top (rows 1-4); existing rows 1-8 shift down to 5-12
$ws.Rows("1:4").Insert()

# Give the new header rows (1-3) the same bordered / wrapped /
# top-left-aligned look as the rest of the table (style index 1 in the
# original workbook), by copying the format from an already-styled cell.
# Row 4 is intentionally left completely untouched - it stays blank with
# no cells at all.
$ws.Range("A5:D5").Copy() | Out-Null
$ws.Range("A1:D3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New title / metadata rows
$ws.Range("A1").Value = "MASTER PACKAGE"

$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "BOM"
$ws.Range("D2").Value = "MASTER"

$ws.Range("A3").Value = "Placements"

# Row 4 stays completely blank (no row element at all in the target)

# Adjust column widths for columns A and B
$ws.Columns(1).ColumnWidth = 26
$ws.Columns(2).ColumnWidth = 21.6
